$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Cells.Item(38, 8).Value = 4777.8
$ws.Cells.Item(38, 9).Value = 2711.1667
$ws.Cells.Item(38, 10).Value = 6155.5557
$ws.Cells.Item(38, 11).Value = 8133.500100000001
$ws.Cells.Item(38, 12).Value = 18466.6671
$ws.Cells.Item(38, 13).Value = -7761.500100000001
$ws.Cells.Item(38, 14).Value = -19210.6671
# Row 69
$ws.Cells.Item(69, 8).Value = 4791.3
$ws.Cells.Item(69, 9).Value = 4637.6665
$ws.Cells.Item(69, 10).Value = 4857.143
$ws.Cells.Item(69, 11).Value = 13912.9995
$ws.Cells.Item(69, 12).Value = 14571.429
$ws.Cells.Item(69, 13).Value = -13038.9995
$ws.Cells.Item(69, 14).Value = -16319.429
# Row 72
$ws.Cells.Item(72, 8).Value = 4791.3
$ws.Cells.Item(72, 9).Value = 4637.6665
$ws.Cells.Item(72, 10).Value = 4857.143
$ws.Cells.Item(72, 11).Value = 41738.9985
$ws.Cells.Item(72, 12).Value = 43714.287
$ws.Cells.Item(72, 13).Value = -37370.9985
$ws.Cells.Item(72, 14).Value = -52450.287
# Row 80
$ws.Cells.Item(80, 8).Value = 1050
$ws.Cells.Item(80, 9).Value = 1104.8235
$ws.Cells.Item(80, 11).Value = 3314.4705
$ws.Cells.Item(80, 13).Value = -2316.4705
# Row 83
$ws.Cells.Item(83, 8).Value = 1050
$ws.Cells.Item(83, 9).Value = 1104.8235
$ws.Cells.Item(83, 11).Value = 9943.4115
$ws.Cells.Item(83, 13).Value = -4951.4115
# Row 92
$ws.Cells.Item(92, 8).Value = 924.5217
$ws.Cells.Item(92, 9).Value = 838.82355
$ws.Cells.Item(92, 10).Value = 1167.3334
$ws.Cells.Item(92, 11).Value = 838.82355
$ws.Cells.Item(92, 12).Value = 1167.3334
$ws.Cells.Item(92, 13).Value = 409.17645
$ws.Cells.Item(92, 14).Value = -3663.3334
# Row 107
$ws.Cells.Item(107, 8).Value = 1977.4
$ws.Cells.Item(107, 9).Value = 1893.5834
$ws.Cells.Item(107, 10).Value = 2103.125
$ws.Cells.Item(107, 11).Value = 1893.5834
$ws.Cells.Item(107, 12).Value = 2103.125
$ws.Cells.Item(107, 13).Value = 26.41660000000002
$ws.Cells.Item(107, 14).Value = -5943.125
# Row 125
$ws.Cells.Item(125, 8).Value = 6160
$ws.Cells.Item(125, 9).Value = 5700
$ws.Cells.Item(125, 11).Value = 51300
$ws.Cells.Item(125, 13).Value = -48840
# Row 138
$ws.Cells.Item(138, 8).Value = 21741642
$ws.Cells.Item(138, 10).Value = 45458724
$ws.Cells.Item(138, 12).Value = 136376172
$ws.Cells.Item(138, 14).Value = -136386452

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Cells.Item(74, 8).Value = 33690
$ws.Cells.Item(74, 9).Value = 42211.137
$ws.Cells.Item(74, 11).Value = 42211.137
$ws.Cells.Item(74, 13).Value = -41337.137
# Row 77
$ws.Cells.Item(77, 8).Value = 33690
$ws.Cells.Item(77, 9).Value = 42211.137
$ws.Cells.Item(77, 11).Value = 211055.685
$ws.Cells.Item(77, 13).Value = -206687.685
# Row 86
$ws.Cells.Item(86, 8).Value = 99999.5
$ws.Cells.Item(86, 9).Value = 99999.5
$ws.Cells.Item(86, 11).Value = 99999.5
$ws.Cells.Item(86, 13).Value = -98813.5
# Row 89
$ws.Cells.Item(89, 8).Value = 99999.5
$ws.Cells.Item(89, 9).Value = 99999.5
$ws.Cells.Item(89, 11).Value = 299998.5
$ws.Cells.Item(89, 13).Value = -294070.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 1802
$ws.Cells.Item(20, 9).Value = 2106.6667
$ws.Cells.Item(20, 11).Value = 2106.6667
$ws.Cells.Item(20, 13).Value = -1859.6667
# Row 26
$ws.Cells.Item(26, 8).Value = 13168.1
$ws.Cells.Item(26, 9).Value = 13168.1
$ws.Cells.Item(26, 11).Value = 13168.1
$ws.Cells.Item(26, 13).Value = -12876.1
# Row 86
$ws.Cells.Item(86, 8).Value = 16762.7
$ws.Cells.Item(86, 9).Value = 10770
$ws.Cells.Item(86, 10).Value = 24087.111
$ws.Cells.Item(86, 11).Value = 10770
$ws.Cells.Item(86, 12).Value = 24087.111
$ws.Cells.Item(86, 13).Value = -9647
$ws.Cells.Item(86, 14).Value = -26333.111
# Row 89
$ws.Cells.Item(89, 8).Value = 16762.7
$ws.Cells.Item(89, 9).Value = 10770
$ws.Cells.Item(89, 10).Value = 24087.111
$ws.Cells.Item(89, 11).Value = 53850
$ws.Cells.Item(89, 12).Value = 120435.555
$ws.Cells.Item(89, 13).Value = -48234
$ws.Cells.Item(89, 14).Value = -131667.555
# Row 141
$ws.Cells.Item(141, 8).Value = 20000
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 38
$ws.Cells.Item(38, 8).Value = 72193
$ws.Cells.Item(38, 9).Value = 166784
$ws.Cells.Item(38, 10).Value = 1249.75
$ws.Cells.Item(38, 11).Value = 166784
$ws.Cells.Item(38, 12).Value = 1249.75
$ws.Cells.Item(38, 13).Value = -166407
$ws.Cells.Item(38, 14).Value = -2003.75
# Row 46
$ws.Cells.Item(46, 8).Value = 72193
$ws.Cells.Item(46, 9).Value = 166784
$ws.Cells.Item(46, 10).Value = 1249.75
$ws.Cells.Item(46, 11).Value = 166784
$ws.Cells.Item(46, 12).Value = 1249.75
$ws.Cells.Item(46, 13).Value = -166573
$ws.Cells.Item(46, 14).Value = -1671.75
# Row 62
$ws.Cells.Item(62, 8).Value = 3690.3215
$ws.Cells.Item(62, 10).Value = 3658.0625
$ws.Cells.Item(62, 12).Value = 3658.0625
$ws.Cells.Item(62, 14).Value = -4906.0625
# Row 65
$ws.Cells.Item(65, 8).Value = 3690.3215
$ws.Cells.Item(65, 10).Value = 3658.0625
$ws.Cells.Item(65, 12).Value = 18290.3125
$ws.Cells.Item(65, 14).Value = -24530.3125

$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Cells.Item(13, 8).Value = 752.5
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 752.5
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 2257.5
$ws.Cells.Item(13, 13).ClearContents()
$ws.Cells.Item(13, 14).Value = -2593.5
# Row 15
$ws.Cells.Item(15, 8).Value = 78.375
$ws.Cells.Item(15, 9).Value = 37.833332
$ws.Cells.Item(15, 10).Value = 200
$ws.Cells.Item(15, 11).Value = 113.499996
$ws.Cells.Item(15, 12).Value = 600
$ws.Cells.Item(15, 13).Value = 26.500004
$ws.Cells.Item(15, 14).Value = -880
# Row 16
$ws.Cells.Item(16, 8).Value = 323.75
$ws.Cells.Item(16, 9).Value = 181.66667
$ws.Cells.Item(16, 10).Value = 750
$ws.Cells.Item(16, 11).Value = 545.00001
$ws.Cells.Item(16, 12).Value = 2250
$ws.Cells.Item(16, 13).Value = -372.00001
$ws.Cells.Item(16, 14).Value = -2596
# Row 22
$ws.Cells.Item(22, 8).Value = 3333.3333
$ws.Cells.Item(22, 9).Value = 3333.3333
$ws.Cells.Item(22, 11).Value = 9999.999899999999
$ws.Cells.Item(22, 13).Value = -9830.999899999999
# Row 27
$ws.Cells.Item(27, 8).Value = 3333.3333
$ws.Cells.Item(27, 9).Value = 3333.3333
$ws.Cells.Item(27, 11).Value = 9999.999899999999
$ws.Cells.Item(27, 13).Value = -9897.999899999999
# Row 75
$ws.Cells.Item(75, 8).Value = 1207.5
$ws.Cells.Item(75, 10).Value = 1815
$ws.Cells.Item(75, 12).Value = 5445
$ws.Cells.Item(75, 14).Value = -7441
# Row 78
$ws.Cells.Item(78, 8).Value = 1207.5
$ws.Cells.Item(78, 10).Value = 1815
$ws.Cells.Item(78, 12).Value = 16335
$ws.Cells.Item(78, 14).Value = -26319
# Row 113
$ws.Cells.Item(113, 8).Value = 1332
$ws.Cells.Item(113, 9).Value = 538.3333
$ws.Cells.Item(113, 11).Value = 1614.9999
$ws.Cells.Item(113, 13).Value = 555.0001
# Row 122
$ws.Cells.Item(122, 8).Value = 649
$ws.Cells.Item(122, 10).Value = 649
$ws.Cells.Item(122, 12).Value = 5841
$ws.Cells.Item(122, 14).Value = -10741
# Row 131
$ws.Cells.Item(131, 8).Value = 998.46875
$ws.Cells.Item(131, 10).Value = 1198.6111
$ws.Cells.Item(131, 12).Value = 3595.8333
$ws.Cells.Item(131, 14).Value = -13675.8333
# Row 137
$ws.Cells.Item(137, 8).Value = 3630
$ws.Cells.Item(137, 10).Value = 3793.3333
$ws.Cells.Item(137, 12).Value = 11379.9999
$ws.Cells.Item(137, 14).Value = -21579.9999

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 1114.48
$ws.Cells.Item(16, 9).Value = 1037.4783
$ws.Cells.Item(16, 11).Value = 1037.4783
$ws.Cells.Item(16, 13).Value = -867.4783
# Row 82
$ws.Cells.Item(82, 8).Value = 3192.6191
$ws.Cells.Item(82, 9).Value = 3134.3845
$ws.Cells.Item(82, 10).Value = 3287.25
$ws.Cells.Item(82, 11).Value = 3134.3845
$ws.Cells.Item(82, 12).Value = 3287.25
$ws.Cells.Item(82, 13).Value = -2773.3845
$ws.Cells.Item(82, 14).Value = -4009.25
# Row 85
$ws.Cells.Item(85, 8).Value = 3192.6191
$ws.Cells.Item(85, 9).Value = 3134.3845
$ws.Cells.Item(85, 10).Value = 3287.25
$ws.Cells.Item(85, 11).Value = 3134.3845
$ws.Cells.Item(85, 12).Value = 3287.25
$ws.Cells.Item(85, 13).Value = -1886.3845
$ws.Cells.Item(85, 14).Value = -5783.25
# Row 122
$ws.Cells.Item(122, 8).Value = 3609.6155
$ws.Cells.Item(122, 9).Value = 3111.818
$ws.Cells.Item(122, 10).Value = 6347.5
$ws.Cells.Item(122, 11).Value = 9335.454000000002
$ws.Cells.Item(122, 12).Value = 19042.5
$ws.Cells.Item(122, 13).Value = -6885.454000000002
$ws.Cells.Item(122, 14).Value = -23942.5

$ws = $wb.Worksheets.Item("WVR")
# Row 56
$ws.Cells.Item(56, 8).Value = 40163.332
$ws.Cells.Item(56, 10).Value = 36990
$ws.Cells.Item(56, 12).Value = 36990
$ws.Cells.Item(56, 14).Value = -38418
# Row 81
$ws.Cells.Item(81, 8).Value = 2620
$ws.Cells.Item(81, 9).Value = 1293
$ws.Cells.Item(81, 11).Value = 2586
$ws.Cells.Item(81, 13).Value = -1525
# Row 82
$ws.Cells.Item(82, 8).Value = 31533.334
$ws.Cells.Item(82, 10).Value = 40000
$ws.Cells.Item(82, 12).Value = 40000
$ws.Cells.Item(82, 14).Value = -40766
# Row 84
$ws.Cells.Item(84, 8).Value = 2620
$ws.Cells.Item(84, 9).Value = 1293
$ws.Cells.Item(84, 11).Value = 12930
$ws.Cells.Item(84, 13).Value = -7626
# Row 85
$ws.Cells.Item(85, 8).Value = 31533.334
$ws.Cells.Item(85, 10).Value = 40000
$ws.Cells.Item(85, 12).Value = 40000
$ws.Cells.Item(85, 14).Value = -42652
# Row 88
$ws.Cells.Item(88, 8).Value = 58999
$ws.Cells.Item(88, 9).Value = 58999
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 58999
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).Value = -58593
$ws.Cells.Item(88, 14).ClearContents()
# Row 91
$ws.Cells.Item(91, 8).Value = 58999
$ws.Cells.Item(91, 9).Value = 58999
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 11).Value = 58999
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 13).Value = -57595
$ws.Cells.Item(91, 14).ClearContents()
# Row 113
$ws.Cells.Item(113, 8).Value = 1281.8572
$ws.Cells.Item(113, 9).Value = 1380.5
$ws.Cells.Item(113, 11).Value = 4141.5
$ws.Cells.Item(113, 13).Value = -1971.5
# Row 126
$ws.Cells.Item(126, 8).Value = 1815
$ws.Cells.Item(126, 10).Value = 2990
$ws.Cells.Item(126, 12).Value = 8970
$ws.Cells.Item(126, 14).Value = -13910
